$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.6008600000000001
$ws.Range("M2").Value = 2.724001666666667
$ws.Range("N2").Value = 8.172005
$ws.Range("O2").Value = 0.04635500474236593
$ws.Range("P2").Value = 0.04635500474236593
$ws.Range("Q2").Value = 1.636743641433333
$ws.Range("R2").Value = 14.7306927729
$ws.Range("S2").Value = 0.04635500474236593
$ws.Range("T2").Value = 0.04635500474236593

# Row 3
$ws.Range("G3").Value = 0.6008600000000001
$ws.Range("O3").Value = 0.6912512390256352
$ws.Range("P3").Value = 0.6912512390256351
$ws.Range("Q3").Value = 24.40731214237334
$ws.Range("S3").Value = 0.6912512390256352
$ws.Range("T3").Value = 0.6912512390256351

# Row 4
$ws.Range("G4").Value = 0.6008600000000001
$ws.Range("M4").Value = 15.419285
$ws.Range("N4").Value = 46.257855
$ws.Range("O4").Value = 0.2623937562319988
$ws.Range("P4").Value = 0.2623937562319988
$ws.Range("Q4").Value = 9.264831585100001
$ws.Range("R4").Value = 83.38348426590001
$ws.Range("S4").Value = 0.2623937562319988
$ws.Range("T4").Value = 0.2623937562319988

$wb.Save()
